$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A (ranking) and column B (vulnerability name) per row
# to match the new order described in the diff. Column C (Soma_Total)
# values stay attached to their original row/text and are left untouched.

$ws.Range("A2").Value = 8
$ws.Range("B2").Value = "Reentrancy"

$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "Integer Overflow and Underflow"

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Transaction Order Dependence"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Assert Violation"

$ws.Range("A6").Value = 6
$ws.Range("B6").Value = "Unchecked Call Return Value"

$ws.Range("A7").Value = 0
$ws.Range("B7").Value = "Unprotected Ether Withdrawal"

$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "Timestamp Dependence"

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Weak Sources of Randomness from Chain Attributes"

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Requirement Violation"

$ws.Range("A11").Value = 11
$ws.Range("B11").Value = "Delegatecall to Untrusted Callee"

$ws.Range("A12").Value = 2
$ws.Range("B12").Value = "Authorization through tx.origin"

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "DoS with Failed Call"

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Write to Arbitrary Storage Location"

$ws.Range("A15").Value = 10
$ws.Range("B15").Value = "Unprotected SELFDESTRUCT Instruction"
